$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add a collapsed "_GoBack" bookmark right after the first paragraph's
#    text ("...YouTube account all in one.") and before its paragraph mark.
#    The engine mis-resolves a bookmark collapsed exactly at a paragraph's
#    trailing boundary, so we temporarily append a marker character, add the
#    bookmark next to it (no longer at the paragraph boundary), then remove
#    the marker again.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1EndMinus1 = $p1.Range.End - 1
$markerPoint = $d.Range($p1EndMinus1, $p1EndMinus1)
$markerPoint.InsertAfter("Z")

$full = $d.Content.Text
$markerIdx = $full.IndexOf("one.Z")
$bmPos = $markerIdx + 4
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$zRange = $d.Range($bmPos, $bmPos + 1)
$zRange.Delete()

# ---------------------------------------------------------------------------
# 2) Merge the four runs "<space>", "The ", "following",
#    " example shows how simple it is to pull from your own site once you've"
#    into a single run, while leaving the following run
#    (" uploaded your XML file...") untouched/separate. A direct Range.Text
#    edit on the whole paragraph merges everything (including the following
#    run) into one run, so we rebuild just this paragraph via InsertXML with
#    the exact run layout we want.
# ---------------------------------------------------------------------------
$targetPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.StartsWith([char]9 + "If you have content")) {
        $targetPara = $para
        break
    }
}
$targetRange = $targetPara.Range

$rssParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t>If you have content that you want to be able to update frequently&#8212;more often than you would want to submit an update to the Window Store&#8212;then you can use your new Google site to upload an XML file, and it&#8217;ll work with the RSS feeds section of your app just like YouTube, etc.</w:t></w:r><w:r><w:t xml:space="preserve"> The following example shows how simple it is to pull from your own site once you&#8217;ve</w:t></w:r><w:r><w:t xml:space="preserve"> uploaded your XML file to your site. (There is a link at the bottom of your Google page with an upload file button.)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$targetRange.InsertXML($rssParaXml)

# ---------------------------------------------------------------------------
# 3) Remove the "_GoBack" bookmark that originally sat after "...you publish
#    the app." (it has been relocated to the end of the first paragraph).
# ---------------------------------------------------------------------------
$oldBookmark = $d.Bookmarks("_GoBack")
if ($oldBookmark.Start -ne $bmPos) {
    $oldBookmark.Delete()
}

# ---------------------------------------------------------------------------
# 4) Merge the three runs "This tutorial is for writing your own RSS feed
#    for your website. ", "If you click ... make things easy.", " " into a
#    single run ending in a single trailing space, keeping the following
#    hyperlink untouched. A plain Range.Text assignment is safe here because
#    the run immediately after the edited range is a <w:hyperlink>, not a
#    plain run, so the engine's run-merge pass does not reach past it.
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$pStart = $full.IndexOf("This tutorial is for")
$pEnd = $full.IndexOf("http://www.w3schools.com")
$helpRange = $d.Range($pStart, $pEnd)

$newHelpText = "This tutorial is for writing your own RSS feed for your website. " + `
    "If you click " + [char]0x201C + "Try it yourself," + [char]0x201D + " then you" + [char]0x2019 + "ll be brought to a place where you can play with the code and see how it will appear. " + `
    "Please note that long descriptions don" + [char]0x2019 + "t seem to go over well with this mini-editor. " + `
    "Just leave that part out for your tests but test everything else about your RSS feed. Use their formatting to make things easy. "

# Force a genuine content change (identical text is treated as a no-op and
# skips run normalization) by appending a temporary marker, then stripping
# it back off.
$helpRange.Text = $newHelpText + "Z"

$full2 = $d.Content.Text
$zIdx = $full2.IndexOf("easy. Z")
$zPos = $zIdx + 6
$helpZRange = $d.Range($zPos, $zPos + 1)
$helpZRange.Delete()

# ---------------------------------------------------------------------------
# 5) Bump the header text's font size to 14pt (sz/szCs = 28) for the
#    paragraph mark and both runs in the header.
# ---------------------------------------------------------------------------
$hdr = $d.Sections(1).Headers(1)
$hdrRange = $hdr.Range
$hdrRange.Font.Size = 14
$hdrRange.Font.SizeBi = 14
